$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns D, H, I, J, K, L, M, P
# (the underlying data rows were reshuffled/updated to reflect the
# weekly consolidation described in the commit message)
$rows = @{
    2  = @{ D = 44497; H = "Sin especificar";  I = "Primera"; J = 160; K = 5000; L = 6000; M = 5500; P = 5500 }
    3  = @{ D = 44263; H = "Sin especificar";  I = "Primera"; J = 100; K = 7000; L = 8000; M = 7500; P = 7500 }
    4  = @{ D = 44259; H = "Sin especificar";  I = "Primera"; J = 80;  K = 4000; L = 4500; M = 4250; P = 4250 }
    5  = @{ D = 44309; H = "Sin especificar";  I = "Primera"; J = 50;  K = 8000; L = 9000; M = 8500; P = 8500 }
    6  = @{ D = 44559; H = "Americana (o)";    I = "Primera"; J = 100; K = 5000; L = 6000; M = 5500; P = 5500 }
    7  = @{ D = 44371; H = "Sin especificar";  I = "Primera"; J = 80;  K = 7000; L = 8000; M = 7375; P = 7375 }
    8  = @{ D = 44253; H = "Americana (o)";    I = "Segunda"; J = 100; K = 4000; L = 4500; M = 4250; P = 4250 }
    9  = @{ D = 44410; H = "Sin especificar";  I = "Primera"; J = 100; K = 5500; L = 6000; M = 5750; P = 5750 }
    10 = @{ D = 44636; H = "Americana (o)";    I = "Primera"; J = 60;  K = 8000; L = 9000; M = 8500; P = 8500 }
    11 = @{ D = 44414; H = "Sin especificar";  I = "Primera"; J = 100; K = 6000; L = 7000; M = 6500; P = 6500 }
    12 = @{ D = 44575; H = "Sin especificar";  I = "Primera"; J = 160; K = 6500; L = 7000; M = 6750; P = 6750 }
    13 = @{ D = 44539; H = "Americana (o)";    I = "Primera"; J = 160; K = 6500; L = 7000; M = 6750; P = 6750 }
    14 = @{ D = 44281; H = "Sin especificar";  I = "Primera"; J = 100; K = 5000; L = 6000; M = 5500; P = 5500 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("H$r").Value = $vals.H
    $ws.Range("I$r").Value = $vals.I
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("P$r").Value = $vals.P
}
